$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Datos actualizados" timestamp text (A1)
$ws.Range("A1").Value = "Datos actualizados a 25 de Mayo de 2020 a las 05:35"

# Update Brasil row (row 5) stats
$ws.Range("B5").Value = 365213
$ws.Range("C5").Value = 1595
$ws.Range("E5").Value = 192556
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 22746

# Insert a new row for Haiti right after "Guinea Ecuatorial" (row 112) and before
# "Niger" (row 113), which pushes Niger / Republica de Chipre / Costa Rica / Zambia
# down by one row. The old Haiti row (previously right before Paraguay) is then
# removed, since it moved up to this new position.
$ws.Rows("113:113").Insert()

# Bring the old Haiti row's content (country label + values) up into the newly
# inserted row.
$ws.Range("A118:H118").Copy()
$ws.Range("A113:H113").PasteSpecial()

# Remove the now-duplicated old Haiti row, closing the gap so all subsequent
# rows (Paraguay, etc.) shift back up to their original row numbers.
$ws.Rows("118:118").Delete()

# Update Haiti's stats (row 113) to the new reported values.
$ws.Range("B113").Value = 958
$ws.Range("C113").Value = 93
$ws.Range("D113").Value = 22
$ws.Range("E113").Value = 909
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 27

# Update Mongolia row (row 162) stats
$ws.Range("D162").Value = 33
$ws.Range("E162").Value = 108
